$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cat")
$ws.Activate()

# EnumValues text: "Yes,No,Maybe" -> "Yes,No,0"
$ws.Range("J13").Value = "Yes,No,0"

# Swap the Required (K) / Nullable (L) flags on rows 7 and 9.
# Using Copy/PasteSpecial (instead of re-typing "true"/"false") keeps the
# cells stored as plain text values (matching the file's existing
# convention) instead of Excel auto-converting typed "true"/"false" into
# native boolean cells.
$scratch = $ws.Range("Z1")

$ws.Range("K7").Copy() | Out-Null
$scratch.PasteSpecial(-4104) | Out-Null
$ws.Range("L7").Copy() | Out-Null
$ws.Range("K7").PasteSpecial(-4104) | Out-Null
$scratch.Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4104) | Out-Null
$scratch.ClearContents() | Out-Null

$ws.Range("K9").Copy() | Out-Null
$scratch.PasteSpecial(-4104) | Out-Null
$ws.Range("L9").Copy() | Out-Null
$ws.Range("K9").PasteSpecial(-4104) | Out-Null
$scratch.Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4104) | Out-Null
$scratch.ClearContents() | Out-Null

# Row 13: Required (K13) flips false -> true; Nullable (L13) already "true".
$ws.Range("L13").Copy() | Out-Null
$ws.Range("K13").PasteSpecial(-4104) | Out-Null

# Move the view's top-left cell and the active selection.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 2
$aw.ScrollRow = 1
$ws.Range("J14").Select() | Out-Null
